## Applies the fixes described in the commit:
## "Added Abstract, FFT subchapter, fixed errors found by ***"
##
## The substantive, content-level change on the "NN time efficiency" sheet
## is that the "Structure" column text for the two rows describing the
## 129-channel / 384x128 classifier runs was wrong: it said only the
## preprocessor ran, but it should say preprocessor *and* postprocessor
## ran (matching the already-correct text used elsewhere in the sheet).
## Also, the previously selected cell on the sheet moves from G14 to F14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foaie1")

$fixedText = "None - only preprocessor and postprocessor ran; returned hardcoded zeros as probabilities"

$ws.Range("C9").Value = $fixedText
$ws.Range("C13").Value = $fixedText

# Excel records an explicit (default) row height on the two edited rows
# once their content is touched.
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15

# Update the active cell/selection recorded in the sheet view.
$ws.Range("F14").Select()

# Record an explicit (portrait) page setup for the sheet.
$ws.PageSetup.Orientation = 1
